# Applies the data refresh for the cryptos list (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.270.29"
$ws.Range("E2").Value = "  +6.89%  "
$ws.Range("D3").Value = "3.682.20"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.43%  "
$ws.Range("D7").Value = "3.680.14"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +4.54%  "
$ws.Range("E10").Value = "  +7.82%  "
$ws.Range("E11").Value = "  +3.48%  "
$ws.Range("E12").Value = "  +6.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.61%  "
$ws.Range("E14").Value = "  +6.31%  "
$ws.Range("D15").Value = "4.298.57"
$ws.Range("E15").Value = "  +19.34%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.684.98"
$ws.Range("E16").Value = "  +19.27%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "71.253.79"
$ws.Range("E17").Value = "  +6.96%  "
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "517.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.18%  "
$ws.Range("E22").Value = "  +17.93%  "
$ws.Range("E23").Value = "  +8.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.14%  "
$ws.Range("E26").Value = "  +8.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.75%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +12.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +18.20%  "
$ws.Range("E34").Value = "  +4.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.342"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.28%  "
$ws.Range("E39").Value = "  +8.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "46.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.77%  "
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("D43").Value = "3.188.74"
$ws.Range("E43").Value = "  +14.87%  "
$ws.Range("E44").Value = "  +6.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "400.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.92%  "
$ws.Range("E47").Value = "  +6.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("E51").Value = "  +11.89%  "
